$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("header")
$ws.Range("A1").Value = "test"
